# Final edit for backlog, sprint and release plan

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 12 (O2): update user-story text and bump the estimate 9 -> 10
$t.Cell(12, 2).Range.Text = "As a user, I want to be able to view the site in my preferred language"
$t.Cell(12, 3).Range.Text = "10"

# Row 14 (O4): replace the user-story text entirely and bump the estimate 3 -> 5
$t.Cell(14, 2).Range.Text = "As a user I want to be able to get access to more variety of information on available city services so that I do not have to look for those information elsewhere."
$t.Cell(14, 3).Range.Text = "5"

# Row 15 (O6 -> O5): renumber the item id
$t.Cell(15, 1).Range.Text = "O5"

# Row 16 (Total row): right-align the "Total" label and bump the grand total 45 -> 48
$t.Cell(16, 2).Range.ParagraphFormat.Alignment = 2
$t.Cell(16, 3).Range.Text = "48"
